# The commit "Fix go demo example" corrects a typo in the "meta" sheet's
# demo row for the map-delimiter example: the value should be "&=" not "=&".
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("meta")

# Fix the typo: B8 changes from "=&" to "&=".
# A leading apostrophe keeps the cell as quoted text (quotePrefix style)
# instead of Excel resetting the cell's formatting when the value is set.
$ws2.Range("B8").Formula = "'&="

# Restore the per-sheet selections that were recorded the last time each
# sheet was active.
$ws2.Activate()
$ws2.Range("C10").Select() | Out-Null

$ws1.Activate()
$ws1.Range("E16").Select() | Out-Null
